$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "to throw away (～を)"
$ws.Range("A3").Value = "to begin (～を)"
$ws.Range("A4").Value = "to drive (～を)"
$ws.Range("A5").Value = "to do laundry (～を)"
$ws.Range("A6").Value = "to clean (～を)"
$ws.Range("A71").Value = "skillful; good at (～が)"
$ws.Range("A72").Value = "clumsy; poor at (～が)"
$ws.Range("A74").Value = "to wash (～を)"
$ws.Range("A76").Value = "to need (～が)"
$ws.Range("A80").Value = "to cut (～を)"
$ws.Range("A81").Value = "to make (～を)"
$ws.Range("A84").Value = "to take (a thing) (～を)"

Write-Output "Done"
